# Certification Project Checklist - add a new checklist row for
# "new user registration" test case (URL of the target site under test).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()

# New row 27: just a URL dropped into column B (no Sr No / Status value
# supplied for this row, matching rows that sit outside the formatted
# table range).
$ws.Range("B27").Value = "http://newtours.demoaut.com/"

# Scroll the view down so row 27 is visible near the top of the window,
# and leave the newly entered cell selected.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select() | Out-Null
